$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert header for the new column I, preserving the original "height(mm)" text
# (do this before renaming E1 so the shared-string slot for "height(mm)" stays in use)
$ws.Range("I1").Value = "height(mm)"

# Rename the old header in E1
$ws.Range("E1").Value = "height(mm)111"

# Correct a data entry in E18 (530 -> 415)
$ws.Range("E18").Value = 415

# Add the new "height(mm)" (adjusted) column I = E - 10 for every data row
for ($r = 2; $r -le 115; $r++) {
    $ws.Range("I$r").Formula = "=E$r-10"
}

# Move the active selection to K9, matching the post-edit cursor position
$ws.Range("K9").Select() | Out-Null
